$d = $word.ActiveDocument

$d.Content.Find.Execute("37+7=44", $true, $false, $false, $false, $false, $true, 1, $false, "26+37=63", 2) | Out-Null
$d.Content.Find.Execute("91-34=57", $true, $false, $false, $false, $false, $true, 1, $false, "19+29=48", 2) | Out-Null
$d.Content.Find.Execute("61-38=23", $true, $false, $false, $false, $false, $true, 1, $false, "24+58=82", 2) | Out-Null
$d.Content.Find.Execute("47+35=82", $true, $false, $false, $false, $false, $true, 1, $false, "45+37=82", 2) | Out-Null
$d.Content.Find.Execute("65-46=19", $true, $false, $false, $false, $false, $true, 1, $false, "77+19=96", 2) | Out-Null
$d.Content.Find.Execute("17+39=56", $true, $false, $false, $false, $false, $true, 1, $false, "97-88=9", 2) | Out-Null
$d.Content.Find.Execute("61-39=22", $true, $false, $false, $false, $false, $true, 1, $false, "64-36=28", 2) | Out-Null
$d.Content.Find.Execute("55+9=64", $true, $false, $false, $false, $false, $true, 1, $false, "80-75=5", 2) | Out-Null
$d.Content.Find.Execute("91-52=39", $true, $false, $false, $false, $false, $true, 1, $false, "63+9=72", 2) | Out-Null
$d.Content.Find.Execute("49+23=72", $true, $false, $false, $false, $false, $true, 1, $false, "66+6=72", 2) | Out-Null
$d.Content.Find.Execute("65-7=58", $true, $false, $false, $false, $false, $true, 1, $false, "28+58=86", 2) | Out-Null
$d.Content.Find.Execute("40-14=26", $true, $false, $false, $false, $false, $true, 1, $false, "34+7=41", 2) | Out-Null
$d.Content.Find.Execute("43-9=34", $true, $false, $false, $false, $false, $true, 1, $false, "92-48=44", 2) | Out-Null
$d.Content.Find.Execute("31-3=28", $true, $false, $false, $false, $false, $true, 1, $false, "67+26=93", 2) | Out-Null
$d.Content.Find.Execute("27+7=34", $true, $false, $false, $false, $false, $true, 1, $false, "27+68=95", 2) | Out-Null
$d.Content.Find.Execute("17+44=61", $true, $false, $false, $false, $false, $true, 1, $false, "71-58=13", 2) | Out-Null
$d.Content.Find.Execute("4+68=72", $true, $false, $false, $false, $false, $true, 1, $false, "77+5=82", 2) | Out-Null
$d.Content.Find.Execute("55-19=36", $true, $false, $false, $false, $false, $true, 1, $false, "30-29=1", 2) | Out-Null
$d.Content.Find.Execute("60-47=13", $true, $false, $false, $false, $false, $true, 1, $false, "14+68=82", 2) | Out-Null
$d.Content.Find.Execute("8+86=94", $true, $false, $false, $false, $false, $true, 1, $false, "48+49=97", 2) | Out-Null
$d.Content.Find.Execute("23+58=81", $true, $false, $false, $false, $false, $true, 1, $false, "61-26=35", 2) | Out-Null
$d.Content.Find.Execute("64-18=46", $true, $false, $false, $false, $false, $true, 1, $false, "52+29=81", 2) | Out-Null
$d.Content.Find.Execute("82-44=38", $true, $false, $false, $false, $false, $true, 1, $false, "15+56=71", 2) | Out-Null
$d.Content.Find.Execute("13-8=5", $true, $false, $false, $false, $false, $true, 1, $false, "69+9=78", 2) | Out-Null
$d.Content.Find.Execute("91-29=62", $true, $false, $false, $false, $false, $true, 1, $false, "38+48=86", 2) | Out-Null
$d.Content.Find.Execute("77-68=9", $true, $false, $false, $false, $false, $true, 1, $false, "81-42=39", 2) | Out-Null
$d.Content.Find.Execute("42-18=24", $true, $false, $false, $false, $false, $true, 1, $false, "91-25=66", 2) | Out-Null
$d.Content.Find.Execute("71-15=56", $true, $false, $false, $false, $false, $true, 1, $false, "28+4=32", 2) | Out-Null
$d.Content.Find.Execute("54-48=6", $true, $false, $false, $false, $false, $true, 1, $false, "28+45=73", 2) | Out-Null
$d.Content.Find.Execute("91-12=79", $true, $false, $false, $false, $false, $true, 1, $false, "29+25=54", 2) | Out-Null
$d.Content.Find.Execute("87-69=18", $true, $false, $false, $false, $false, $true, 1, $false, "21-14=7", 2) | Out-Null
$d.Content.Find.Execute("29+8=37", $true, $false, $false, $false, $false, $true, 1, $false, "92-54=38", 2) | Out-Null
$d.Content.Find.Execute("34+28=62", $true, $false, $false, $false, $false, $true, 1, $false, "82-8=74", 2) | Out-Null
$d.Content.Find.Execute("56+36=92", $true, $false, $false, $false, $false, $true, 1, $false, "51-37=14", 2) | Out-Null
$d.Content.Find.Execute("52-8=44", $true, $false, $false, $false, $false, $true, 1, $false, "7+7=14", 2) | Out-Null
$d.Content.Find.Execute("93-24=69", $true, $false, $false, $false, $false, $true, 1, $false, "62-53=9", 2) | Out-Null
$d.Content.Find.Execute("18+39=57", $true, $false, $false, $false, $false, $true, 1, $false, "26+49=75", 2) | Out-Null
$d.Content.Find.Execute("64+7=71", $true, $false, $false, $false, $false, $true, 1, $false, "30-1=29", 2) | Out-Null
$d.Content.Find.Execute("45+7=52", $true, $false, $false, $false, $false, $true, 1, $false, "59+17=76", 2) | Out-Null
$d.Content.Find.Execute("22+59=81", $true, $false, $false, $false, $false, $true, 1, $false, "84-36=48", 2) | Out-Null
$d.Content.Find.Execute("74-29=45", $true, $false, $false, $false, $false, $true, 1, $false, "74-9=65", 2) | Out-Null
$d.Content.Find.Execute("51-48=3", $true, $false, $false, $false, $false, $true, 1, $false, "42-17=25", 2) | Out-Null
$d.Content.Find.Execute("27+65=92", $true, $false, $false, $false, $false, $true, 1, $false, "70-26=44", 2) | Out-Null
$d.Content.Find.Execute("42+39=81", $true, $false, $false, $false, $false, $true, 1, $false, "3+18=21", 2) | Out-Null
$d.Content.Find.Execute("23+18=41", $true, $false, $false, $false, $false, $true, 1, $false, "82-78=4", 2) | Out-Null
$d.Content.Find.Execute("19+56=75", $true, $false, $false, $false, $false, $true, 1, $false, "18+17=35", 2) | Out-Null
$d.Content.Find.Execute("80-2=78", $true, $false, $false, $false, $false, $true, 1, $false, "60-11=49", 2) | Out-Null
$d.Content.Find.Execute("89+5=94", $true, $false, $false, $false, $false, $true, 1, $false, "87-58=29", 2) | Out-Null
$d.Content.Find.Execute("13-6=7", $true, $false, $false, $false, $false, $true, 1, $false, "92-47=45", 2) | Out-Null
$d.Content.Find.Execute("56-39=17", $true, $false, $false, $false, $false, $true, 1, $false, "73-68=5", 2) | Out-Null
$d.Content.Find.Execute("48+25=73", $true, $false, $false, $false, $false, $true, 1, $false, "64-28=36", 2) | Out-Null
$d.Content.Find.Execute("68+26=94", $true, $false, $false, $false, $false, $true, 1, $false, "73-47=26", 2) | Out-Null
$d.Content.Find.Execute("15+38=53", $true, $false, $false, $false, $false, $true, 1, $false, "48-19=29", 2) | Out-Null
$d.Content.Find.Execute("44-26=18", $true, $false, $false, $false, $false, $true, 1, $false, "81-36=45", 2) | Out-Null
$d.Content.Find.Execute("65-58=7", $true, $false, $false, $false, $false, $true, 1, $false, "84-49=35", 2) | Out-Null
$d.Content.Find.Execute("7+75=82", $true, $false, $false, $false, $false, $true, 1, $false, "53-37=16", 2) | Out-Null
$d.Content.Find.Execute("85-49=36", $true, $false, $false, $false, $false, $true, 1, $false, "90-82=8", 2) | Out-Null
$d.Content.Find.Execute("55-29=26", $true, $false, $false, $false, $false, $true, 1, $false, "72-58=14", 2) | Out-Null
$d.Content.Find.Execute("55+16=71", $true, $false, $false, $false, $false, $true, 1, $false, "26+16=42", 2) | Out-Null
$d.Content.Find.Execute("33-7=26", $true, $false, $false, $false, $false, $true, 1, $false, "38+55=93", 2) | Out-Null
$d.Content.Find.Execute("16+19=35", $true, $false, $false, $false, $false, $true, 1, $false, "82-19=63", 2) | Out-Null
$d.Content.Find.Execute("46-29=17", $true, $false, $false, $false, $false, $true, 1, $false, "41-34=7", 2) | Out-Null
$d.Content.Find.Execute("81-3=78", $true, $false, $false, $false, $false, $true, 1, $false, "61-28=33", 2) | Out-Null
$d.Content.Find.Execute("17+15=32", $true, $false, $false, $false, $false, $true, 1, $false, "28+29=57", 2) | Out-Null
$d.Content.Find.Execute("95-86=9", $true, $false, $false, $false, $false, $true, 1, $false, "5+26=31", 2) | Out-Null
$d.Content.Find.Execute("27+38=65", $true, $false, $false, $false, $false, $true, 1, $false, "95-19=76", 2) | Out-Null
$d.Content.Find.Execute("28+48=76", $true, $false, $false, $false, $false, $true, 1, $false, "18+59=77", 2) | Out-Null
$d.Content.Find.Execute("28+54=82", $true, $false, $false, $false, $false, $true, 1, $false, "71-13=58", 2) | Out-Null
$d.Content.Find.Execute("13+49=62", $true, $false, $false, $false, $false, $true, 1, $false, "3+78=81", 2) | Out-Null
$d.Content.Find.Execute("44-35=9", $true, $false, $false, $false, $false, $true, 1, $false, "6+85=91", 2) | Out-Null
$d.Content.Find.Execute("22+39=61", $true, $false, $false, $false, $false, $true, 1, $false, "50-39=11", 2) | Out-Null
$d.Content.Find.Execute("78-9=69", $true, $false, $false, $false, $false, $true, 1, $false, "4+57=61", 2) | Out-Null
$d.Content.Find.Execute("86+7=93", $true, $false, $false, $false, $false, $true, 1, $false, "23-6=17", 2) | Out-Null
$d.Content.Find.Execute("95-88=7", $true, $false, $false, $false, $false, $true, 1, $false, "27+37=64", 2) | Out-Null
$d.Content.Find.Execute("40-9=31", $true, $false, $false, $false, $false, $true, 1, $false, "35+6=41", 2) | Out-Null
$d.Content.Find.Execute("72-66=6", $true, $false, $false, $false, $false, $true, 1, $false, "14-8=6", 2) | Out-Null
$d.Content.Find.Execute("15+18=33", $true, $false, $false, $false, $false, $true, 1, $false, "18+23=41", 2) | Out-Null
$d.Content.Find.Execute("29+66=95", $true, $false, $false, $false, $false, $true, 1, $false, "87+9=96", 2) | Out-Null
$d.Content.Find.Execute("28+25=53", $true, $false, $false, $false, $false, $true, 1, $false, "17+27=44", 2) | Out-Null
$d.Content.Find.Execute("74-5=69", $true, $false, $false, $false, $false, $true, 1, $false, "17+38=55", 2) | Out-Null
$d.Content.Find.Execute("28+44=72", $true, $false, $false, $false, $false, $true, 1, $false, "56+16=72", 2) | Out-Null
$d.Content.Find.Execute("8+79=87", $true, $false, $false, $false, $false, $true, 1, $false, "60-28=32", 2) | Out-Null
$d.Content.Find.Execute("34+39=73", $true, $false, $false, $false, $false, $true, 1, $false, "96-9=87", 2) | Out-Null
$d.Content.Find.Execute("14+59=73", $true, $false, $false, $false, $false, $true, 1, $false, "16+56=72", 2) | Out-Null
$d.Content.Find.Execute("69+14=83", $true, $false, $false, $false, $false, $true, 1, $false, "56+39=95", 2) | Out-Null
$d.Content.Find.Execute("66-17=49", $true, $false, $false, $false, $false, $true, 1, $false, "19+69=88", 2) | Out-Null
$d.Content.Find.Execute("33-26=7", $true, $false, $false, $false, $false, $true, 1, $false, "58+29=87", 2) | Out-Null
$d.Content.Find.Execute("81-23=58", $true, $false, $false, $false, $false, $true, 1, $false, "57-8=49", 2) | Out-Null
$d.Content.Find.Execute("65+28=93", $true, $false, $false, $false, $false, $true, 1, $false, "81-54=27", 2) | Out-Null
$d.Content.Find.Execute("70-3=67", $true, $false, $false, $false, $false, $true, 1, $false, "50-14=36", 2) | Out-Null
$d.Content.Find.Execute("67+28=95", $true, $false, $false, $false, $false, $true, 1, $false, "33-25=8", 2) | Out-Null
$d.Content.Find.Execute("79+6=85", $true, $false, $false, $false, $false, $true, 1, $false, "70-62=8", 2) | Out-Null
$d.Content.Find.Execute("11-8=3", $true, $false, $false, $false, $false, $true, 1, $false, "25+59=84", 2) | Out-Null
$d.Content.Find.Execute("41-27=14", $true, $false, $false, $false, $false, $true, 1, $false, "51-32=19", 2) | Out-Null
$d.Content.Find.Execute("8+36=44", $true, $false, $false, $false, $false, $true, 1, $false, "48+49=97", 2) | Out-Null
$d.Content.Find.Execute("17+16=33", $true, $false, $false, $false, $false, $true, 1, $false, "85-27=58", 2) | Out-Null
$d.Content.Find.Execute("56-48=8", $true, $false, $false, $false, $false, $true, 1, $false, "57+38=95", 2) | Out-Null
$d.Content.Find.Execute("68+8=76", $true, $false, $false, $false, $false, $true, 1, $false, "60-22=38", 2) | Out-Null
$d.Content.Find.Execute("53-4=49", $true, $false, $false, $false, $false, $true, 1, $false, "90-22=68", 2) | Out-Null
$d.Content.Find.Execute("21-2=19", $true, $false, $false, $false, $false, $true, 1, $false, "91-57=34", 2) | Out-Null
